$d = $word.ActiveDocument

# Each entry: text to locate (old line) and the text to put in its place.
# We locate the range via Find (no in-place replace, wdReplaceNone = 0) and
# then assign Range.Text directly -- this avoids Word's "replace as you
# type" smart-quote autocorrect that Find.Execute's ReplaceWith parameter
# would otherwise trigger, so the literal straight quotes are preserved.
$replacements = @(
    @{old = '"football betting odds" → /sport/betting/uk/football-odds.htm'; new = '"football betting odds" → /sport/betting/uk/football-betting-sites.htm'},
    @{old = '"welcome bonus offers" → /sport/betting/uk/welcome-bonuses.htm'; new = '"welcome bonus offers" → /sport/betting/uk/betting-offers.htm'},
    @{old = '"live betting guide" → /sport/betting/live-betting-guide.htm'; new = '"live betting guide" → /sport/betting/guides/live-betting.htm'},
    @{old = '"odds comparison" → /sport/betting/odds-comparison.htm'; new = '"odds comparison" → /sport/betting/uk/index.htm'},
    @{old = '"casino bonuses UK" → /sport/betting/uk/casino-bonuses.htm'; new = '"casino bonuses UK" → /sport/betting/uk/betting-offers.htm'},
    @{old = '"payment methods guide" → /sport/betting/payment-methods.htm'; new = '"payment methods guide" → /sport/betting/uk/index.htm'},
    @{old = '"how to place a bet" → /sport/betting/how-to-place-bet.htm'; new = '"how to place a bet" → /sport/betting/guides/first-bet.htm'},
    @{old = '"betting strategies" → /sport/betting/strategies.htm'; new = '"betting strategies" → /sport/betting/guides/parlay.htm'},
    @{old = '"sports betting glossary" → /sport/betting/glossary.htm'; new = '"sports betting glossary" → /sport/betting/guides/betting-glossary.htm'}
)

foreach ($rep in $replacements) {
    $range = $d.Content
    $found = $range.Find.Execute($rep.old, $true, $true, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
    if ($found) {
        $range.Text = $rep.new
    } else {
        throw "Could not find text: $($rep.old)"
    }
}

$d.Save()
